$wb = $excel.ActiveWorkbook

# Sheet1: Ringkasan Umum
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2,1).Value = 39
$ws1.Cells.Item(2,2).Value = 10
$ws1.Cells.Item(2,3).Value = 10
$ws1.Cells.Item(2,4).Value = 6

# Sheet2: Akses per Jam
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2,1).Value = 3
$ws2.Cells.Item(2,2).Value = 1
$ws2.Cells.Item(3,1).Value = 4
$ws2.Cells.Item(3,2).Value = 1
$ws2.Cells.Item(4,1).Value = 5
$ws2.Cells.Item(4,2).Value = 4
$ws2.Cells.Item(5,1).Value = 7
$ws2.Cells.Item(5,2).Value = 1
$ws2.Cells.Item(6,1).Value = 8
$ws2.Cells.Item(6,2).Value = 4
$ws2.Cells.Item(7,1).Value = 9
$ws2.Cells.Item(7,2).Value = 3
$ws2.Cells.Item(8,1).Value = 10
$ws2.Cells.Item(8,2).Value = 6
$ws2.Cells.Item(9,1).Value = 11
$ws2.Cells.Item(9,2).Value = 1
$ws2.Cells.Item(10,1).Value = 12
$ws2.Cells.Item(10,2).Value = 5
$ws2.Cells.Item(11,1).Value = 13
$ws2.Cells.Item(11,2).Value = 1
$ws2.Cells.Item(12,1).Value = 14
$ws2.Cells.Item(12,2).Value = 1
$ws2.Cells.Item(13,1).Value = 15
$ws2.Cells.Item(13,2).Value = 1
$ws2.Cells.Item(14,1).Value = 16
$ws2.Cells.Item(14,2).Value = 1
$ws2.Cells.Item(15,1).Value = 17
$ws2.Cells.Item(15,2).Value = 3
$ws2.Cells.Item(16,1).Value = 18
$ws2.Cells.Item(16,2).Value = 2
$ws2.Cells.Item(17,1).Value = 19
$ws2.Cells.Item(17,2).Value = 2
$ws2.Cells.Item(18,1).Value = 20
$ws2.Cells.Item(18,2).Value = 1
$ws2.Cells.Item(19,1).Value = 21
$ws2.Cells.Item(19,2).Value = 1

# Sheet3: Akses per Tanggal
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2,1).Value = 45815
$ws3.Cells.Item(2,1).NumberFormat = "YYYY-MM-DD"
$ws3.Cells.Item(2,2).Value = 4
$ws3.Cells.Item(3,1).Value = 45816
$ws3.Cells.Item(3,1).NumberFormat = "YYYY-MM-DD"
$ws3.Cells.Item(3,2).Value = 3
$ws3.Cells.Item(4,1).Value = 45817
$ws3.Cells.Item(4,1).NumberFormat = "YYYY-MM-DD"
$ws3.Cells.Item(4,2).Value = 6
$ws3.Cells.Item(5,1).Value = 45819
$ws3.Cells.Item(5,1).NumberFormat = "YYYY-MM-DD"
$ws3.Cells.Item(5,2).Value = 3
$ws3.Cells.Item(6,1).Value = 45820
$ws3.Cells.Item(6,1).NumberFormat = "YYYY-MM-DD"
$ws3.Cells.Item(6,2).Value = 5
$ws3.Cells.Item(7,1).Value = 45821
$ws3.Cells.Item(7,1).NumberFormat = "YYYY-MM-DD"
$ws3.Cells.Item(7,2).Value = 1
$ws3.Cells.Item(8,1).Value = 45822
$ws3.Cells.Item(8,1).NumberFormat = "YYYY-MM-DD"
$ws3.Cells.Item(8,2).Value = 1
$ws3.Cells.Item(9,1).Value = 45823
$ws3.Cells.Item(9,1).NumberFormat = "YYYY-MM-DD"
$ws3.Cells.Item(9,2).Value = 5
$ws3.Cells.Item(10,1).Value = 45825
$ws3.Cells.Item(10,1).NumberFormat = "YYYY-MM-DD"
$ws3.Cells.Item(10,2).Value = 1
$ws3.Cells.Item(11,1).Value = 45826
$ws3.Cells.Item(11,1).NumberFormat = "YYYY-MM-DD"
$ws3.Cells.Item(11,2).Value = 8
$ws3.Cells.Item(12,1).Value = 45827
$ws3.Cells.Item(12,1).NumberFormat = "YYYY-MM-DD"
$ws3.Cells.Item(12,2).Value = 2
# Remove extra row (old row7 data no longer needed since it moved to row12, but we have 11 rows now vs old 7, so no deletion needed, only growth)

# Sheet4: Top 10 Holder
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2,1).Value = "Ahmad Nur Hidayat"
$ws4.Cells.Item(2,2).Value = 5
$ws4.Cells.Item(3,1).Value = "Master"
$ws4.Cells.Item(3,2).Value = 5
$ws4.Cells.Item(4,1).Value = "Ingrid Vera Mour"
$ws4.Cells.Item(4,2).Value = 5
$ws4.Cells.Item(5,1).Value = "Nafis"
$ws4.Cells.Item(5,2).Value = 5
$ws4.Cells.Item(6,1).Value = "Housekeeping"
$ws4.Cells.Item(6,2).Value = 5
$ws4.Cells.Item(7,1).Value = "Shift Midle"
$ws4.Cells.Item(7,2).Value = 4
$ws4.Cells.Item(8,1).Value = "Triswantoro"
$ws4.Cells.Item(8,2).Value = 3
$ws4.Cells.Item(9,1).Value = "Metia Shanti Wulanda"
$ws4.Cells.Item(9,2).Value = 3
$ws4.Cells.Item(10,1).Value = "Orbani"
$ws4.Cells.Item(10,2).Value = 2
$ws4.Cells.Item(11,1).Value = "Lucky Hendrawan"
$ws4.Cells.Item(11,2).Value = 1

# Sheet5: Penggunaan Kartu
$ws5 = $wb.Worksheets.Item(5)
$ws5.Cells.Item(2,1).Value = 60
$ws5.Cells.Item(2,2).Value = "Master"
$ws5.Cells.Item(2,3).Value = "Master Card"
$ws5.Cells.Item(2,4).Value = 5
$ws5.Cells.Item(3,1).Value = 10875
$ws5.Cells.Item(3,2).Value = "Housekeeping"
$ws5.Cells.Item(3,3).Value = "HSKP Card"
$ws5.Cells.Item(3,4).Value = 5
$ws5.Cells.Item(4,1).Value = 13313
$ws5.Cells.Item(4,2).Value = "Ingrid Vera Mour"
$ws5.Cells.Item(4,3).Value = "Guest Card"
$ws5.Cells.Item(4,4).Value = 5
$ws5.Cells.Item(5,1).Value = 13274
$ws5.Cells.Item(5,2).Value = "Nafis"
$ws5.Cells.Item(5,3).Value = "Guest Card"
$ws5.Cells.Item(5,4).Value = 5
$ws5.Cells.Item(6,1).Value = 13423
$ws5.Cells.Item(6,2).Value = "Ahmad Nur Hidayat"
$ws5.Cells.Item(6,3).Value = "Guest Card"
$ws5.Cells.Item(6,4).Value = 5
$ws5.Cells.Item(7,1).Value = 7516
$ws5.Cells.Item(7,2).Value = "Shift Midle"
$ws5.Cells.Item(7,3).Value = "HSKP Card"
$ws5.Cells.Item(7,4).Value = 4
$ws5.Cells.Item(8,1).Value = 13386
$ws5.Cells.Item(8,2).Value = "Metia Shanti Wulanda"
$ws5.Cells.Item(8,3).Value = "Guest Card"
$ws5.Cells.Item(8,4).Value = 3
$ws5.Cells.Item(9,1).Value = 13344
$ws5.Cells.Item(9,2).Value = "Triswantoro"
$ws5.Cells.Item(9,3).Value = "Guest Card"
$ws5.Cells.Item(9,4).Value = 3
$ws5.Cells.Item(10,1).Value = 13333
$ws5.Cells.Item(10,2).Value = "Orbani"
$ws5.Cells.Item(10,3).Value = "Guest Card"
$ws5.Cells.Item(10,4).Value = 2
$ws5.Cells.Item(11,1).Value = 13429
$ws5.Cells.Item(11,2).Value = "Lucky Hendrawan"
$ws5.Cells.Item(11,3).Value = "Guest Card"
$ws5.Cells.Item(11,4).Value = 1

# Sheet6: Log Dini Hari - delete rows 9, 8 (bottom-up), then overwrite rows 2-7
$ws6 = $wb.Worksheets.Item(6)
$ws6.Rows.Item(9).Delete()
$ws6.Rows.Item(8).Delete()
$ws6.Cells.Item(2,1).Value = 87
$ws6.Cells.Item(2,2).Value = 45826.24097222222
$ws6.Cells.Item(2,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws6.Cells.Item(2,3).ClearContents()
$ws6.Cells.Item(2,4).Value = "Guest Card"
$ws6.Cells.Item(2,5).Value = 13423
$ws6.Cells.Item(2,6).Value = "Ahmad Nur Hidayat"
$ws6.Cells.Item(2,7).ClearContents()
$ws6.Cells.Item(2,8).Value = "vino"
$ws6.Cells.Item(2,9).Value = "18/06/2025 05:40:40"
$ws6.Cells.Item(2,10).ClearContents()
$ws6.Cells.Item(2,11).ClearContents()
$ws6.Cells.Item(2,12).Value = 45826
$ws6.Cells.Item(2,12).NumberFormat = "YYYY-MM-DD"
$ws6.Cells.Item(2,13).Value = 5
$ws6.Cells.Item(2,14).Value = "2025-06-18 05:47:00"
$ws6.Cells.Item(3,1).Value = 95
$ws6.Cells.Item(3,2).Value = 45826.23611111111
$ws6.Cells.Item(3,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws6.Cells.Item(3,3).ClearContents()
$ws6.Cells.Item(3,4).Value = "Guest Card"
$ws6.Cells.Item(3,5).Value = 13423
$ws6.Cells.Item(3,6).Value = "Ahmad Nur Hidayat"
$ws6.Cells.Item(3,7).ClearContents()
$ws6.Cells.Item(3,8).Value = "vino"
$ws6.Cells.Item(3,9).Value = "18/06/2025 05:40:40"
$ws6.Cells.Item(3,10).ClearContents()
$ws6.Cells.Item(3,11).ClearContents()
$ws6.Cells.Item(3,12).Value = 45826
$ws6.Cells.Item(3,12).NumberFormat = "YYYY-MM-DD"
$ws6.Cells.Item(3,13).Value = 5
$ws6.Cells.Item(3,14).Value = "2025-06-18 05:40:00"
$ws6.Cells.Item(4,1).Value = 108
$ws6.Cells.Item(4,2).Value = 45826.23194444444
$ws6.Cells.Item(4,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws6.Cells.Item(4,3).ClearContents()
$ws6.Cells.Item(4,4).Value = "Guest Card"
$ws6.Cells.Item(4,5).Value = 13423
$ws6.Cells.Item(4,6).Value = "Ahmad Nur Hidayat"
$ws6.Cells.Item(4,7).ClearContents()
$ws6.Cells.Item(4,8).Value = "vino"
$ws6.Cells.Item(4,9).Value = "18/06/2025 05:40:40"
$ws6.Cells.Item(4,10).ClearContents()
$ws6.Cells.Item(4,11).ClearContents()
$ws6.Cells.Item(4,12).Value = 45826
$ws6.Cells.Item(4,12).NumberFormat = "YYYY-MM-DD"
$ws6.Cells.Item(4,13).Value = 5
$ws6.Cells.Item(4,14).Value = "2025-06-18 05:34:00"
$ws6.Cells.Item(5,1).Value = 173
$ws6.Cells.Item(5,2).Value = 45825.14722222222
$ws6.Cells.Item(5,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws6.Cells.Item(5,3).ClearContents()
$ws6.Cells.Item(5,4).Value = "Master Card"
$ws6.Cells.Item(5,5).Value = 60
$ws6.Cells.Item(5,6).Value = "Master"
$ws6.Cells.Item(5,7).Value = "front office"
$ws6.Cells.Item(5,8).Value = "hotel"
$ws6.Cells.Item(5,9).Value = "07/01/2023 22:19:13"
$ws6.Cells.Item(5,10).Value = "hotel"
$ws6.Cells.Item(5,11).Value = "07/01/2025 06:33:52"
$ws6.Cells.Item(5,12).Value = 45825
$ws6.Cells.Item(5,12).NumberFormat = "YYYY-MM-DD"
$ws6.Cells.Item(5,13).Value = 3
$ws6.Cells.Item(5,14).Value = "2025-06-17 03:32:00"
$ws6.Cells.Item(6,1).Value = 313
$ws6.Cells.Item(6,2).Value = 45820.22638888889
$ws6.Cells.Item(6,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws6.Cells.Item(6,3).ClearContents()
$ws6.Cells.Item(6,4).Value = "Master Card"
$ws6.Cells.Item(6,5).Value = 60
$ws6.Cells.Item(6,6).Value = "Master"
$ws6.Cells.Item(6,7).Value = "front office"
$ws6.Cells.Item(6,8).Value = "hotel"
$ws6.Cells.Item(6,9).Value = "07/01/2023 22:19:13"
$ws6.Cells.Item(6,10).Value = "hotel"
$ws6.Cells.Item(6,11).Value = "07/01/2025 06:33:52"
$ws6.Cells.Item(6,12).Value = 45820
$ws6.Cells.Item(6,12).NumberFormat = "YYYY-MM-DD"
$ws6.Cells.Item(6,13).Value = 5
$ws6.Cells.Item(6,14).Value = "2025-06-12 05:26:00"
$ws6.Cells.Item(7,1).Value = 391
$ws6.Cells.Item(7,2).Value = 45817.18263888889
$ws6.Cells.Item(7,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws6.Cells.Item(7,3).ClearContents()
$ws6.Cells.Item(7,4).Value = "Guest Card"
$ws6.Cells.Item(7,5).Value = 13313
$ws6.Cells.Item(7,6).Value = "Ingrid Vera Mour"
$ws6.Cells.Item(7,7).ClearContents()
$ws6.Cells.Item(7,8).Value = "eko"
$ws6.Cells.Item(7,9).Value = "09/06/2025 04:26:48"
$ws6.Cells.Item(7,10).ClearContents()
$ws6.Cells.Item(7,11).ClearContents()
$ws6.Cells.Item(7,12).Value = 45817
$ws6.Cells.Item(7,12).NumberFormat = "YYYY-MM-DD"
$ws6.Cells.Item(7,13).Value = 4
$ws6.Cells.Item(7,14).Value = "2025-06-09 04:23:00"
